$d = $word.ActiveDocument

# Locate the bibliography paragraph ("COMMELLI, Leo. Lista de comandos
# uteis do GIT, ..."). The new "Link do projeto no GitHub: ..." paragraph
# must be inserted immediately before it (i.e. right after the existing
# blank paragraph that follows the "git status" line).
$searchRange = $d.Content
$null = $searchRange.Find.Execute("COMMELLI, Leo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetPara = $searchRange.Paragraphs(1)
$targetRange = $targetPara.Range
$targetRange.Collapse(1)
$targetRange.InsertParagraphBefore()

# Re-resolve the bibliography paragraph (the previous handles may be stale
# now that the document has been mutated) and grab the paragraph that was
# just created right before it.
$searchRange2 = $d.Content
$null = $searchRange2.Find.Execute("COMMELLI, Leo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$newParaRange = $searchRange2.Paragraphs(1).Previous().Range

# Populate the new (already correctly-formatted, inherited Times New Roman
# / 12pt / justified) paragraph with two runs: the label and the URL.
$xml = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Link do projeto no GitHub: </w:t></w:r>' +
              '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>https://github.com/FelipeCarvalho25/ProgramMobileJF.git</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$newParaRange.InsertXML($xml)
